$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 77
$newRow = 78

# Insert a new row below the last data row, copying the formatting (and values,
# which get overwritten below) of the last data row.
$ws.Range("A" + $lastRow + ":J" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":J" + $newRow).Insert(-4121)  # xlShiftDown
$excel.CutCopyMode = $false

# Fill in the new row's values
$ws.Cells.Item($newRow, 1).Value = 43978
$ws.Cells.Item($newRow, 2).Value = 77210
$ws.Cells.Item($newRow, 3).Value = 631
$ws.Cells.Item($newRow, 4).Value = 1473
$ws.Cells.Item($newRow, 5).Value = 2
$ws.Cells.Item($newRow, 6).Value = 7
$ws.Cells.Item($newRow, 7).Value = 2
$ws.Cells.Item($newRow, 8).Value = 1
$ws.Cells.Item($newRow, 9).Value = 108
$ws.Cells.Item($newRow, 10).Value = 0

# Expand the table range to include the new row
$tbl = $ws.ListObjects.Item("Tabela1")
$tbl.Resize($ws.Range("A1:J" + $newRow))

# Update the selection to mirror the recorded change
$ws.Range("A" + $newRow + ":J" + $newRow).Select() | Out-Null
